# Re-run the UnmetDemand model output generation: the underlying demand
# input for this scenario was zeroed out, so all of the derived dispatch /
# demand / revenue figures in the generated "Output" workbook collapse to
# zero (except Total Profits, which settles at the negative of the fixed
# operating cost that still accrues with no revenue/variable cost).

$wb = $excel.ActiveWorkbook

# --- "Costs and Revenues" --------------------------------------------------
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2:D2").Value = 0        # Total Revenues
$ws.Range("B4:D4").Value = 0        # Total Operation Variable Costs
$ws.Range("B6:D6").Value = -33627.6 # Total Profits

# --- "Connected Households" -------------------------------------------------
$ws = $wb.Worksheets.Item("Connected Households")
$ws.Range("B2:D2").Value = 0        # Consumers
$ws.Range("B3:D3").Value = 0        # Prosumers

# --- "DG Dispatch" ------------------------------------------------------
$ws = $wb.Worksheets.Item("DG Dispatch")
$ws.Range("B2:Y10").Value = 0

# --- "Fed-in Capacity" ---------------------------------------------------
$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("I2:S10").Value = 0

# --- "Yearly demand" ------------------------------------------------------
$ws = $wb.Worksheets.Item("Yearly demand")
$ws.Range("B2:Y10").Value = 0

# --- "Net demand" ---------------------------------------------------------
$ws = $wb.Worksheets.Item("Net demand")
$ws.Range("B2:Y10").Value = 0

# --- "Net surplus" ---------------------------------------------------------
$ws = $wb.Worksheets.Item("Net surplus")
$ws.Range("I2:S10").Value = 0

# --- "Unmet Demand" ---------------------------------------------------------
$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("B2:Y10").Value = 0
